$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.284.15"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "1.888.45"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.96"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4692"
$ws.Range("E7").Value = "  -1.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2845"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06607"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.81"
$ws.Range("E10").Value = "  +11.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07783"
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.91"
$ws.Range("E12").Value = "  -3.24%  "
$ws.Range("D13").Value = "1.881.36"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.118"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6753"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "284.23"
$ws.Range("E16").Value = "  +10.89%  "
$ws.Range("D17").Value = "30.295.30"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9997"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.67"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.407"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007307"
$ws.Range("E22").Value = "  -2.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.184"
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.421"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.24"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.28"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.996"
$ws.Range("E28").Value = "  -3.03%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09752"
$ws.Range("E30").Value = "  -3.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.388"
$ws.Range("E31").Value = "  -7.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.488"
$ws.Range("E32").Value = "  -1.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.150"
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04686"
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7085"
$ws.Range("E35").Value = "  -2.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.099"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.710"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01873"
$ws.Range("E38").Value = "  -2.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.732"
$ws.Range("E39").Value = "  +7.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.532"
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.34"
$ws.Range("E41").Value = "  -2.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.974"
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8695"
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.91"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4194"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("D47").Value = "1.002.62"
$ws.Range("E47").Value = "  +10.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.299"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.188"
$ws.Range("E49").Value = "  +5.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.97"
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1152"
$ws.Range("E51").Value = "  -3.99%  "
